# Generate Report for Handoff
# Adds a new handed-off file (b2ffb6d6-5cbb-4a78-bec5-c7e7534f1bf0.md) as a
# new row to the "Overview", "zh-cn" and "de-de" tables/sheets, mirroring the
# existing 0243ac06-... row.

$wb = $excel.ActiveWorkbook

$commit = "73a7f375d4cf6d9ed2b105818f94c421b468aca1"
$fileBase = "b2ffb6d6-5cbb-4a78-bec5-c7e7534f1bf0"
$ghUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/$fileBase.md"

# ---------------------------------------------------------------------
# Sheet "Overview" (table "Overview") -> new row 3
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item("Overview")
$loOv.ListRows.Add() | Out-Null

$wsOv.Range("A3").Value = "$fileBase.md"
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), $ghUrl, "", "", "e2e\$fileBase.md") | Out-Null
$wsOv.Range("C3").Value = ".md"
$wsOv.Range("E3").Value = "Ready for handoff"
$wsOv.Range("F3").Value = "Ready for handoff"
$wsOv.Range("G3").Value = "2016-09-04 14:43:18"
$wsOv.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "zh-cn" (table "zh-cn") -> new row 3
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item("zh-cn")
$loZh.ListRows.Add() | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $ghUrl, "", "", "$fileBase.md") | Out-Null
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Formula = "'False"
$wsZh.Range("G3").Value = "$fileBase.779220ce3fa2dc58c47bf444b24f0c1b243a8d2f.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-09-04 14:43:13"
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("M3").Formula = "'True"
$wsZh.Range("O3").Formula = "'False"

# ---------------------------------------------------------------------
# Sheet "de-de" (table "de-de") -> new row 3
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item("de-de")
$loDe.ListRows.Add() | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $ghUrl, "", "", "$fileBase.md") | Out-Null
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Formula = "'False"
$wsDe.Range("G3").Value = "$fileBase.779220ce3fa2dc58c47bf444b24f0c1b243a8d2f.de-de.xlf"
$wsDe.Range("H3").Value = "2016-09-04 14:43:18"
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M3").Formula = "'True"
$wsDe.Range("O3").Formula = "'False"

Write-Host "Report row added for $fileBase"
